$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14.71780674912695
$ws.Range("C2").Value = 11.70139857436583
$ws.Range("E2").Value = 16.65366325213244
$ws.Range("F2").Value = 16.86991607391245
$ws.Range("G2").Value = 13.04860781319793
$ws.Range("H2").Value = 10.35190583803758
$ws.Range("I2").Value = 13.84225927835601
$ws.Range("O2").Value = 13.63216782526221

$ws.Range("B3").Value = 13.83054310157664
$ws.Range("C3").Value = 11.12281386937704
$ws.Range("E3").Value = 15.69970040925352
$ws.Range("F3").Value = 15.89584955866815
$ws.Range("G3").Value = 13.18392910761913
$ws.Range("H3").Value = 10.42636047001692
$ws.Range("I3").Value = 14.00779444967782
$ws.Range("O3").Value = 13.76448272730525

$ws.Range("B4").Value = 13.25414245937379
$ws.Range("C4").Value = 10.75051581907269
$ws.Range("E4").Value = 15.0880924750741
$ws.Range("F4").Value = 15.26997757108489
$ws.Range("G4").Value = 13.28040293957676
$ws.Range("H4").Value = 10.47480516671016
$ws.Range("I4").Value = 14.11433864705918
$ws.Range("O4").Value = 13.85152333248498

$ws.Range("B5").Value = 13.01138256781458
$ws.Range("C5").Value = 10.59464622650531
$ws.Range("E5").Value = 14.83261628750293
$ws.Range("F5").Value = 15.008197319934
$ws.Range("G5").Value = 13.32301202290392
$ws.Range("H5").Value = 10.49523254448105
$ws.Range("I5").Value = 14.15899364093858
$ws.Range("O5").Value = 13.88844263345184

$ws.Range("B6").Value = 12.97059975214662
$ws.Range("C6").Value = 10.56851753142524
$ws.Range("E6").Value = 14.78982649855366
$ws.Range("F6").Value = 14.96433081551589
$ws.Range("G6").Value = 13.33028415876488
$ws.Range("H6").Value = 10.49866589979107
$ws.Range("I6").Value = 14.16648341968669
$ws.Range("O6").Value = 13.89466035879916

$ws.Range("B7").Value = 13.25090026836361
$ws.Range("C7").Value = 10.74843034091548
$ws.Range("E7").Value = 15.08467191494239
$ws.Range("F7").Value = 15.26647399323133
$ws.Range("G7").Value = 13.28096433292298
$ws.Range("H7").Value = 10.47507788148405
$ws.Range("I7").Value = 14.11493586462363
$ws.Range("O7").Value = 13.85201538111991

$ws.Range("B8").Value = 14.41848103643941
$ws.Range("C8").Value = 11.50551172751504
$ws.Range("E8").Value = 16.3302493538257
$ws.Range("F8").Value = 16.5399640634477
$ws.Range("G8").Value = 13.09244617496611
$ws.Range("H8").Value = 10.37701101282191
$ws.Range("I8").Value = 13.89832019423237
$ws.Range("O8").Value = 13.67658123433209

$ws.Range("B9").Value = 16.47436723651634
$ws.Range("C9").Value = 12.85024392944021
$ws.Range("E9").Value = 18.69378146740698
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 12.83204773975139
$ws.Range("H9").Value = 10.20638192326153
$ws.Range("I9").Value = 13.51227205749852
$ws.Range("O9").Value = 13.37896201063405

$ws.Range("B10").Value = 17.84342886432501
$ws.Range("C10").Value = 13.74797551823884
$ws.Range("E10").Value = 20.3538630016205
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 12.71136816850399
$ws.Range("H10").Value = 10.09427249213954
$ws.Range("I10").Value = 13.25199559172523
$ws.Range("O10").Value = 13.18915422001561

$ws.Range("B11").Value = 18.43027371351761
$ws.Range("C11").Value = 14.13603732319945
$ws.Range("E11").Value = 21.06632429473285
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 12.672593360517
$ws.Range("H11").Value = 10.04615765325466
$ws.Range("I11").Value = 13.13860427766371
$ws.Range("O11").Value = 13.10919900102106

$ws.Range("B12").Value = 18.64733510500001
$ws.Range("C12").Value = 14.2800174355442
$ws.Range("E12").Value = 21.33000886585616
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 12.66028724084899
$ws.Range("H12").Value = 10.0283535871441
$ws.Range("I12").Value = 13.09638215316278
$ws.Range("O12").Value = 13.07985208697082

$ws.Range("B13").Value = 18.60081659988498
$ws.Range("C13").Value = 14.24914152150564
$ws.Range("E13").Value = 21.27349084249297
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 12.66283098826836
$ws.Range("H13").Value = 10.03216949013489
$ws.Range("I13").Value = 13.10544362600105
$ws.Range("O13").Value = 13.08613090131038

$ws.Range("B14").Value = 18.4482349365256
$ws.Range("C14").Value = 14.14794242867473
$ws.Range("E14").Value = 21.08814002520047
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 12.6715329723925
$ws.Range("H14").Value = 10.04468455786437
$ws.Range("I14").Value = 13.13511629529753
$ws.Range("O14").Value = 13.10676589833717

$ws.Range("B15").Value = 18.35410197259271
$ws.Range("C15").Value = 14.08556695156986
$ws.Range("E15").Value = 20.97381271584676
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 12.67717439969619
$ws.Range("H15").Value = 10.0524046087538
$ws.Range("I15").Value = 13.15338488901453
$ws.Range("O15").Value = 13.11952694594723

$ws.Range("B16").Value = 17.80435417213255
$ws.Range("C16").Value = 13.72220149693032
$ws.Range("E16").Value = 20.30644527362895
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 12.7142317575196
$ws.Range("H16").Value = 10.09747505175044
$ws.Range("I16").Value = 13.25950647378065
$ws.Range("O16").Value = 13.19450903480515

$ws.Range("B17").Value = 17.45789708209684
$ws.Range("C17").Value = 13.49404626663578
$ws.Range("E17").Value = 19.88611992035678
$ws.Range("F17").Value = 20.20408069597325
$ws.Range("G17").Value = 12.74113805310558
$ws.Range("H17").Value = 10.12586388697585
$ws.Range("I17").Value = 13.32588911366982
$ws.Range("O17").Value = 13.24215271190501

$ws.Range("B18").Value = 17.25523993554082
$ws.Range("C18").Value = 13.36090691594229
$ws.Range("E18").Value = 19.64033705944609
$ws.Range("F18").Value = 19.95656407809801
$ws.Range("G18").Value = 12.75812820907727
$ws.Range("H18").Value = 10.14246379974239
$ws.Range("I18").Value = 13.36454241337279
$ws.Range("O18").Value = 13.27015676605225

$ws.Range("B19").Value = 17.18604236975246
$ws.Range("C19").Value = 13.31550172905581
$ws.Range("E19").Value = 19.55642724761393
$ws.Range("F19").Value = 19.87204792380568
$ws.Range("G19").Value = 12.76413900702364
$ws.Range("H19").Value = 10.14813082352324
$ws.Range("I19").Value = 13.37771089481551
$ws.Range("O19").Value = 13.27974127596995

$ws.Range("B20").Value = 17.49512825500455
$ws.Range("C20").Value = 13.51853189145298
$ws.Range("E20").Value = 19.93128037485902
$ws.Range("F20").Value = 20.24955283636154
$ws.Range("G20").Value = 12.73811665941588
$ws.Range("H20").Value = 10.12281374896752
$ws.Range("I20").Value = 13.31877376959721
$ws.Range("O20").Value = 13.23701871096639

$ws.Range("B21").Value = 18.49319195985865
$ws.Range("C21").Value = 14.17774801008754
$ws.Range("E21").Value = 21.14274758342711
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 12.66891203578207
$ws.Range("H21").Value = 10.04099727936591
$ws.Range("I21").Value = 13.12638129118179
$ws.Range("O21").Value = 13.10067955436983

$ws.Range("B22").Value = 19.11538925927481
$ws.Range("C22").Value = 14.59125382119352
$ws.Range("E22").Value = 21.89891866981113
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 12.63756133782923
$ws.Range("H22").Value = 9.989950764033614
$ws.Range("I22").Value = 13.00481749407803
$ws.Range("O22").Value = 13.01700274412497

$ws.Range("B23").Value = 18.78606191501423
$ws.Range("C23").Value = 14.37215720727207
$ws.Range("E23").Value = 21.49858157779473
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("G23").Value = 12.65300598141578
$ws.Range("H23").Value = 10.01697290980718
$ws.Range("I23").Value = 13.06931750755024
$ws.Range("O23").Value = 13.06116196366763

$ws.Range("B24").Value = 17.47830686399793
$ws.Range("C24").Value = 13.50746807249492
$ws.Range("E24").Value = 19.91087620750696
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 12.73947789700916
$ws.Range("H24").Value = 10.12419184797679
$ws.Range("I24").Value = 13.32198909296829
$ws.Range("O24").Value = 13.23933788390448

$ws.Range("B25").Value = 15.93830940449256
$ws.Range("C25").Value = 12.50200977734833
$ws.Range("E25").Value = 18.04433547837734
$ws.Range("F25").Value = 18.34778573295695
$ws.Range("G25").Value = 12.89036723027032
$ws.Range("H25").Value = 10.25021731566778
$ws.Range("I25").Value = 13.61258852132513
$ws.Range("O25").Value = 13.45445104844271
